# Add data for 2022-08-14 (carjacking-by-neighborhood-by-month)
#
# This updates the "through" date (sheet name + header cell B1) from
# August 05 to August 06, and bumps a handful of per-neighborhood /
# per-month cell counts to reflect the newly-added day's incidents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- rename the sheet and update the rolling "through" header --------------
$ws.Name = "Through 2022-08-06"
$ws.Range("B1").Value = "August 2022 (through August 06)"

# --- helper to add a delta to a cell's current numeric value ---------------
function Bump-Cell([string]$addr, [double]$delta) {
    $cell = $ws.Range($addr)
    $current = $cell.Value()
    if ($null -eq $current) { $current = 0 }
    $cell.Value = $current + $delta
}

# Row 2
Bump-Cell "B2" 1
Bump-Cell "AH2" 1
Bump-Cell "BF2" 1

# Row 3
Bump-Cell "AP3" 1

# Row 4
Bump-Cell "AP4" 2

# Row 5
Bump-Cell "B5" 1
Bump-Cell "J5" 2
Bump-Cell "R5" 1

# Row 6
Bump-Cell "B6" 1
Bump-Cell "J6" 1

# Row 7
Bump-Cell "AX7" 1

# Row 9
Bump-Cell "AP9" 1

# Row 14
Bump-Cell "AP14" 1

# Row 15
Bump-Cell "B15" 1

# Row 20
Bump-Cell "R20" 1
Bump-Cell "AP20" 1

# Row 29
Bump-Cell "J29" 1

# Row 30
Bump-Cell "J30" 1

# Row 32
Bump-Cell "B32" 1

# Row 33
Bump-Cell "J33" 1

# Row 36
Bump-Cell "B36" 1

# Row 37
Bump-Cell "R37" 1

# Row 44
Bump-Cell "R44" 1

# Row 45
Bump-Cell "J45" 1

# Row 46
Bump-Cell "Z46" 1

# Row 57
Bump-Cell "J57" 1

# Row 82
Bump-Cell "J82" 1

# Row 95
Bump-Cell "J95" 1
